$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay text (Excel would otherwise
# auto-convert numeric-looking text like "606.60" -> 606.6). Force text format
# for the whole Price column up front, matching the original inline-string data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: 'Bitcoin'
$ws.Range('D2').Value = '66.558.15'
$ws.Range('E2').Value = '  -0.25%  '

# Row 3: 'Ethereum'
$ws.Range('D3').Value = '3.523.09'
$ws.Range('E3').Value = '  -3.07%  '

# Row 4: 'TetherUSD'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5: 'BNB'
$ws.Range('D5').Value = '606.60'
$ws.Range('E5').Value = '  -0.68%  '

# Row 6: 'Solana'
$ws.Range('D6').Value = '143.84'
$ws.Range('E6').Value = '  -4.18%  '

# Row 7: 'LidoStakedEther'
$ws.Range('D7').Value = '3.522.07'
$ws.Range('E7').Value = '  -3.09%  '

# Row 8: 'USDC'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.21%  '

# Row 9: 'XRP'
$ws.Range('E9').Value = '  +3.43%  '

# Row 10: 'Toncoin'
$ws.Range('D10').Value = '7.72'
$ws.Range('E10').Value = '  -2.93%  '

# Row 11: 'Dogecoin'
$ws.Range('E11').Value = '  -4.96%  '

# Row 12: 'Cardano'
$ws.Range('E12').Value = '  -2.91%  '

# Row 13: 'WrappedliquidstakedEther2.0'
$ws.Range('D13').Value = '4.108.49'
$ws.Range('E13').Value = '  -3.33%  '

# Row 14: 'ShibaInu'
$ws.Range('E14').Value = '  -6.91%  '

# Row 15: 'Avalanche'
$ws.Range('D15').Value = '28.70'
$ws.Range('E15').Value = '  -4.33%  '

# Row 16: 'WrappedEther'
$ws.Range('D16').Value = '3.539.91'
$ws.Range('E16').Value = '  -2.56%  '

# Row 17: 'TRON'
$ws.Range('E17').Value = '  -0.08%  '

# Row 18: 'WrappedBTC'
$ws.Range('D18').Value = '66.368.47'
$ws.Range('E18').Value = '  -0.69%  '

# Row 19: 'Uniswap'
$ws.Range('D19').Value = '10.76'
$ws.Range('E19').Value = '  -7.60%  '

# Row 20: 'Polkadot'
$ws.Range('D20').Value = '6.14'
$ws.Range('E20').Value = '  -3.86%  '

# Row 21: 'Chainlink'
$ws.Range('D21').Value = '14.61'
$ws.Range('E21').Value = '  -3.70%  '

# Row 22: 'BitcoinCash'
$ws.Range('D22').Value = '422.84'
$ws.Range('E22').Value = '  -1.45%  '

# Row 23: 'Polygon'
$ws.Range('D23').Value = '0.590'
$ws.Range('E23').Value = '  -5.27%  '

# Row 24: 'Litecoin'
$ws.Range('D24').Value = '77.03'
$ws.Range('E24').Value = '  -2.41%  '

# Row 25: 'WrappedeETH'
$ws.Range('D25').Value = '3.668.08'
$ws.Range('E25').Value = '  -2.87%  '

# Row 26: 'Dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.01%  '

# Row 27: 'PEPE'
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').Value = '  -8.01%  '

# Row 28: 'PancakeSwap'
$ws.Range('E28').Value = '  -2.63%  '

# Row 29: 'RenderToken'
$ws.Range('D29').Value = '7.84'
$ws.Range('E29').Value = '  -6.80%  '

# Row 30: 'InternetComputer(DFINITY)'
$ws.Range('D30').Value = '8.91'
$ws.Range('E30').Value = '  -6.68%  '

# Row 31: 'Binance-PegBSC-USD'
$ws.Range('E31').Value = '  +0.13%  '

# Row 32: 'RenzoRestakedETH'
$ws.Range('D32').Value = '3.527.17'
$ws.Range('E32').Value = '  -2.89%  '

# Row 33: 'Kaspa'
$ws.Range('D33').Value = '0.154'
$ws.Range('E33').Value = '  -3.03%  '

# Row 34: 'EthereumClassic'
$ws.Range('D34').Value = '24.23'
$ws.Range('E34').Value = '  -5.01%  '

# Row 35: 'USDe'
$ws.Range('E35').Value = '  +0.02%  '

# Row 36: 'Fetch.AI'
$ws.Range('E36').Value = '  -10.26%  '

# Row 37: 'Aptos'
$ws.Range('D37').Value = '7.55'
$ws.Range('E37').Value = '  -4.69%  '

# Row 38: 'ImmutableX'
$ws.Range('E38').Value = '  -5.25%  '

# Row 39: 'Monero'
$ws.Range('D39').Value = '173.07'
$ws.Range('E39').Value = '  -2.10%  '

# Row 40: 'NEARProtocol'
$ws.Range('D40').Value = '5.19'
$ws.Range('E40').Value = '  -9.00%  '

# Row 41: 'Hedera'
$ws.Range('E41').Value = '  -6.46%  '

# Row 42: 'Filecoin'
$ws.Range('E42').Value = '  -5.39%  '

# Row 43: 'Mantle'
$ws.Range('D43').Value = '0.852'
$ws.Range('E43').Value = '  -5.60%  '

# Row 44: 'OKB'
$ws.Range('D44').Value = '45.56'
$ws.Range('E44').Value = '  -0.71%  '

# Row 45: 'Stacks'
$ws.Range('E45').Value = '  -7.03%  '

# Row 46: 'FirstDigitalUSD'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.01%  '

# Row 47: 'dogwifhat'
$ws.Range('E47').Value = '  -8.65%  '

# Row 48: 'Cosmos'
$ws.Range('D48').Value = '7.06'
$ws.Range('E48').Value = '  -2.24%  '

# Row 49: 'EnergySwap' -> 'ONDO'
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '1.12'
$ws.Range('E49').Value = '  -5.79%  '

# Row 50: 'ONDO' -> 'EnergySwap'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '22.90'
$ws.Range('E50').Value = '  -4.96%  '

# Row 51: 'SuiNetwork'
$ws.Range('D51').Value = '0.903'
$ws.Range('E51').Value = '  -6.15%  '
